$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 16 data: column A = index 14 (styled like the other index cells),
# column B = same label as row 15 (HexGrid-60degTilt5degRes),
# columns C:M = averaged intensity values.

# Copy the formatting (bold, bordered, centered) from A15 to A16 first,
# then set the value.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$values = @(
    1.674365088669574,
    1.549275029639975,
    0.6158939211221109,
    1.674365088669574,
    1.228437297018139,
    1.20694505993909,
    0.7380088131999583,
    1.549275029639975,
    1.082584475381043,
    1.378474782025308,
    1.168820868264808
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = 3 + $i  # C = 3
    $ws.Cells.Item(16, $col).Value = $values[$i]
}
